# Apply the scheduled set of text replacements to the document.
# Each "old" value is unique within the document, so a simple
# Find/Replace (non-wildcard, match whole content) per pair is safe.

$d = $word.ActiveDocument

$replacements = @(
    @{ old = "2025-07-08 Tuesday"; new = "2025-07-09 Wednesday" },
    @{ old = "54×17="; new = "83×93=" },
    @{ old = "50×13="; new = "47×91=" },
    @{ old = "11×71="; new = "44×15=" },
    @{ old = "93×69="; new = "72×80=" },
    @{ old = "77×23="; new = "74×42=" },
    @{ old = "48×71="; new = "16×18=" },
    @{ old = "59×50="; new = "18×26=" },
    @{ old = "42×37="; new = "51×42=" },
    @{ old = "65×55="; new = "52×94=" },
    @{ old = "11×81="; new = "54×63=" },
    @{ old = "90×66="; new = "65×31=" },
    @{ old = "40×96="; new = "14×82=" },
    @{ old = "51×92="; new = "91×34=" },
    @{ old = "29×31="; new = "23×34=" },
    @{ old = "62×75="; new = "29×84=" },
    @{ old = "75×51="; new = "26×80=" },
    @{ old = "14×47="; new = "33×77=" },
    @{ old = "71×65="; new = "54×69=" },
    @{ old = "16×84="; new = "61×83=" },
    @{ old = "44×90="; new = "45×69=" },
    @{ old = "13×46="; new = "77×13=" },
    @{ old = "54×29="; new = "45×30=" },
    @{ old = "47×21="; new = "97×29=" },
    @{ old = "91×58="; new = "92×65=" },
    @{ old = "40×44="; new = "42×68=" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
